$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.045.56"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "3.006.90"
$ws.Range("E3").Value = "  +3.02%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.63"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.01"
$ws.Range("E6").Value = "  -3.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.559"
$ws.Range("E7").Value = "  -1.40%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.612"
$ws.Range("E9").Value = "  -2.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.11"
$ws.Range("E10").Value = "  -2.84%  "
$ws.Range("E11").Value = "  +2.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0859"
$ws.Range("E12").Value = "  -3.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.03"
$ws.Range("E13").Value = "  -3.46%  "
$ws.Range("D14").Value = "3.481.35"
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.64"
$ws.Range("E15").Value = "  -3.36%  "
$ws.Range("D16").Value = "3.017.07"
$ws.Range("E16").Value = "  +3.94%  "
$ws.Range("E17").Value = "  +4.00%  "
$ws.Range("D18").Value = "52.130.51"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("E19").Value = "  +3.99%  "
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.58"
$ws.Range("E21").Value = "  -2.80%  "
$ws.Range("D22").Value = "0.0₃0972"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.18"
$ws.Range("E23").Value = "  -2.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.85"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("E25").Value = "  -3.37%  "
$ws.Range("E26").Value = "  -1.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.00"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.45"
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.46"
$ws.Range("E31").Value = "  +6.58%  "
$ws.Range("E32").Value = "  -3.57%  "
$ws.Range("E33").Value = "  -5.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.17"
$ws.Range("E34").Value = "  +15.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.11"
$ws.Range("E35").Value = "  -2.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0437"
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("E38").Value = "  +3.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.86"
$ws.Range("E39").Value = "  +4.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.96"
$ws.Range("E40").Value = "  -2.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.62"
$ws.Range("E41").Value = "  -4.49%  "
$ws.Range("E42").Value = "  -1.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.09"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "124.26"
$ws.Range("E44").Value = "  +2.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.14"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("D46").Value = "2.123.99"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.33"
$ws.Range("E47").Value = "  -3.11%  "
$ws.Range("E48").Value = "  -6.53%  "
$ws.Range("D49").Value = "3.304.74"
$ws.Range("E49").Value = "  +3.06%  "
$ws.Range("E50").Value = "  -2.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0337"
$ws.Range("E51").Value = "  +1.22%  "
